$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(2, 8).Value = 250
$ws.Cells.Item(2, 10).Value = 150
$ws.Cells.Item(2, 12).Value = 150
$ws.Cells.Item(2, 14).Value = -376
$ws.Cells.Item(9, 8).Value = 150
$ws.Cells.Item(9, 9).Value = 0
$ws.Cells.Item(9, 11).Value = 0
$ws.Cells.Item(9, 13).ClearContents()
$ws.Cells.Item(38, 8).Value = 1240.75
$ws.Cells.Item(38, 9).Value = 1007.8
$ws.Cells.Item(38, 10).Value = 1629
$ws.Cells.Item(38, 11).Value = 3023.4
$ws.Cells.Item(38, 12).Value = 4887
$ws.Cells.Item(38, 13).Value = -2651.4
$ws.Cells.Item(38, 14).Value = -5631
$ws.Cells.Item(82, 8).Value = 4116.5
$ws.Cells.Item(82, 9).Value = 4116.5
$ws.Cells.Item(82, 11).Value = 12349.5
$ws.Cells.Item(82, 13).Value = -11943.5
$ws.Cells.Item(85, 8).Value = 4116.5
$ws.Cells.Item(85, 9).Value = 4116.5
$ws.Cells.Item(85, 11).Value = 12349.5
$ws.Cells.Item(85, 13).Value = -10945.5
$ws.Cells.Item(86, 8).Value = 2400.1667
$ws.Cells.Item(86, 9).Value = 2499.25
$ws.Cells.Item(86, 11).Value = 2499.25
$ws.Cells.Item(86, 13).Value = -1376.25
$ws.Cells.Item(89, 8).Value = 2400.1667
$ws.Cells.Item(89, 9).Value = 2499.25
$ws.Cells.Item(89, 11).Value = 12496.25
$ws.Cells.Item(89, 13).Value = -6880.25
$ws.Cells.Item(125, 8).Value = 621.8570999999999
$ws.Cells.Item(125, 9).Value = 465.75
$ws.Cells.Item(125, 10).Value = 830
$ws.Cells.Item(125, 11).Value = 4191.75
$ws.Cells.Item(125, 12).Value = 7470
$ws.Cells.Item(125, 13).Value = -1731.75
$ws.Cells.Item(125, 14).Value = -12390
$ws.Cells.Item(127, 8).Value = 5768
$ws.Cells.Item(127, 9).Value = 5768
$ws.Cells.Item(127, 11).Value = 17304
$ws.Cells.Item(127, 13).Value = -12344
$ws.Cells.Item(132, 8).Value = 1169.48
$ws.Cells.Item(132, 9).Value = 1169.48
$ws.Cells.Item(132, 11).Value = 3508.44
$ws.Cells.Item(132, 13).Value = -978.4400000000001

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(45, 9).Value = 18001368
$ws.Cells.Item(45, 10).Value = 1785.75
$ws.Cells.Item(45, 11).Value = 18001368
$ws.Cells.Item(45, 12).Value = 1785.75
$ws.Cells.Item(45, 13).Value = -18000991
$ws.Cells.Item(45, 14).Value = -2539.75
$ws.Cells.Item(61, 8).Value = 4508.5
$ws.Cells.Item(61, 9).Value = 3114
$ws.Cells.Item(61, 11).Value = 3114
$ws.Cells.Item(61, 13).Value = -2902
$ws.Cells.Item(74, 8).Value = 4899
$ws.Cells.Item(74, 9).Value = 4899
$ws.Cells.Item(74, 11).Value = 4899
$ws.Cells.Item(74, 13).Value = -4025
$ws.Cells.Item(77, 8).Value = 4899
$ws.Cells.Item(77, 9).Value = 4899
$ws.Cells.Item(77, 11).Value = 24495
$ws.Cells.Item(77, 13).Value = -20127
$ws.Cells.Item(136, 8).Value = 4508.5
$ws.Cells.Item(136, 9).Value = 3114
$ws.Cells.Item(136, 11).Value = 9342
$ws.Cells.Item(136, 13).Value = -6792

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(86, 8).Value = 224466.67
$ws.Cells.Item(86, 9).Value = 2500
$ws.Cells.Item(86, 10).Value = 668400
$ws.Cells.Item(86, 11).Value = 2500
$ws.Cells.Item(86, 12).Value = 668400
$ws.Cells.Item(86, 13).Value = -1377
$ws.Cells.Item(86, 14).Value = -670646
$ws.Cells.Item(89, 8).Value = 224466.67
$ws.Cells.Item(89, 9).Value = 2500
$ws.Cells.Item(89, 10).Value = 668400
$ws.Cells.Item(89, 11).Value = 12500
$ws.Cells.Item(89, 12).Value = 3342000
$ws.Cells.Item(89, 13).Value = -6884
$ws.Cells.Item(89, 14).Value = -3353232
$ws.Cells.Item(94, 8).Value = 696.6
$ws.Cells.Item(94, 9).Value = 521.8182
$ws.Cells.Item(94, 11).Value = 521.8182
$ws.Cells.Item(94, 13).Value = -70.81820000000005

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(16, 8).Value = 932.6667
$ws.Cells.Item(16, 10).Value = 0
$ws.Cells.Item(16, 12).Value = 0
$ws.Cells.Item(16, 14).ClearContents()
$ws.Cells.Item(52, 8).Value = 64780
$ws.Cells.Item(52, 10).Value = 64780
$ws.Cells.Item(52, 12).Value = 64780
$ws.Cells.Item(52, 14).Value = -65368
$ws.Cells.Item(58, 8).Value = 0
$ws.Cells.Item(58, 9).Value = 0
$ws.Cells.Item(58, 11).Value = 0
$ws.Cells.Item(58, 13).ClearContents()
$ws.Cells.Item(113, 8).Value = 932.6667
$ws.Cells.Item(113, 10).Value = 0
$ws.Cells.Item(113, 12).Value = 0
$ws.Cells.Item(113, 14).ClearContents()
$ws.Cells.Item(132, 8).Value = 3243.7693
$ws.Cells.Item(132, 9).Value = 1192.8334
$ws.Cells.Item(132, 11).Value = 3578.5002
$ws.Cells.Item(132, 13).Value = -1048.5002
$ws.Cells.Item(134, 8).Value = 3250.9
$ws.Cells.Item(134, 9).Value = 2701.25
$ws.Cells.Item(134, 11).Value = 8103.75
$ws.Cells.Item(134, 13).Value = -5568.75
$ws.Cells.Item(136, 8).Value = 0
$ws.Cells.Item(136, 9).Value = 0
$ws.Cells.Item(136, 11).Value = 0
$ws.Cells.Item(136, 13).ClearContents()
$ws.Cells.Item(138, 8).Value = 24500
$ws.Cells.Item(138, 10).Value = 24500
$ws.Cells.Item(138, 12).Value = 24500
$ws.Cells.Item(138, 14).Value = -34780

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(2, 8).Value = 549.75
$ws.Cells.Item(2, 10).Value = 549.75
$ws.Cells.Item(2, 12).Value = 3298.5
$ws.Cells.Item(2, 14).Value = -3524.5
$ws.Cells.Item(7, 8).Value = 597.5
$ws.Cells.Item(7, 9).Value = 200
$ws.Cells.Item(7, 10).Value = 730
$ws.Cells.Item(7, 11).Value = 600
$ws.Cells.Item(7, 12).Value = 2190
$ws.Cells.Item(7, 13).Value = -488
$ws.Cells.Item(7, 14).Value = -2414
$ws.Cells.Item(82, 8).Value = 6500
$ws.Cells.Item(82, 10).Value = 6500
$ws.Cells.Item(82, 12).Value = 19500
$ws.Cells.Item(82, 14).Value = -20312
$ws.Cells.Item(85, 8).Value = 6500
$ws.Cells.Item(85, 10).Value = 6500
$ws.Cells.Item(85, 12).Value = 19500
$ws.Cells.Item(85, 14).Value = -22308
$ws.Cells.Item(86, 8).Value = 1219.8
$ws.Cells.Item(86, 9).Value = 2250
$ws.Cells.Item(86, 10).Value = 533
$ws.Cells.Item(86, 11).Value = 6750
$ws.Cells.Item(86, 12).Value = 1599
$ws.Cells.Item(86, 13).Value = -5564
$ws.Cells.Item(86, 14).Value = -3971
$ws.Cells.Item(89, 8).Value = 1219.8
$ws.Cells.Item(89, 9).Value = 2250
$ws.Cells.Item(89, 10).Value = 533
$ws.Cells.Item(89, 11).Value = 20250
$ws.Cells.Item(89, 12).Value = 4797
$ws.Cells.Item(89, 13).Value = -14322
$ws.Cells.Item(89, 14).Value = -16653
$ws.Cells.Item(131, 8).Value = 10015045
$ws.Cells.Item(131, 10).Value = 17408.791
$ws.Cells.Item(131, 12).Value = 52226.37300000001
$ws.Cells.Item(131, 14).Value = -62306.37300000001

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(102, 8).Value = 2146.1365
$ws.Cells.Item(102, 9).Value = 2320.6365
$ws.Cells.Item(102, 10).Value = 1971.6364
$ws.Cells.Item(102, 11).Value = 2320.6365
$ws.Cells.Item(102, 12).Value = 1971.6364
$ws.Cells.Item(102, 13).Value = -698.6365000000001
$ws.Cells.Item(102, 14).Value = -5215.6364

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(40, 8).Value = 4471.241
$ws.Cells.Item(40, 9).Value = 2358.35
$ws.Cells.Item(40, 10).Value = 9166.556
$ws.Cells.Item(40, 11).Value = 2358.35
$ws.Cells.Item(40, 12).Value = 9166.556
$ws.Cells.Item(40, 13).Value = -2222.35
$ws.Cells.Item(40, 14).Value = -9438.556
$ws.Cells.Item(93, 8).Value = 944.5
$ws.Cells.Item(93, 9).Value = 813.5625
$ws.Cells.Item(93, 11).Value = 813.5625
$ws.Cells.Item(93, 13).Value = 434.4375
$ws.Cells.Item(136, 8).Value = 3248.7646
$ws.Cells.Item(136, 9).Value = 2113.087
$ws.Cells.Item(136, 11).Value = 6339.261
$ws.Cells.Item(136, 13).Value = -3789.261
$ws.Cells.Item(139, 8).Value = 43933.332
$ws.Cells.Item(139, 9).Value = 40000
$ws.Cells.Item(139, 11).Value = 40000
$ws.Cells.Item(139, 13).Value = -34860
